$d = $word.ActiveDocument

$replacements = @(
    @('42×60=2520', '34×85=2890'),
    @('89×82=7298', '84×31=2604'),
    @('93×78=7254', '91×30=2730'),
    @('19×89=1691', '98×79=7742'),
    @('48×20=960', '59×80=4720'),
    @('35×86=3010', '34×40=1360'),
    @('68×18=1224', '92×16=1472'),
    @('49×65=3185', '88×95=8360'),
    @('88×63=5544', '60×83=4980'),
    @('97×65=6305', '67×84=5628'),
    @('33×30=990', '64×52=3328'),
    @('98×81=7938', '53×29=1537'),
    @('11×97=1067', '38×51=1938'),
    @('98×90=8820', '26×94=2444'),
    @('92×94=8648', '59×96=5664'),
    @('28×49=1372', '42×76=3192'),
    @('54×61=3294', '45×25=1125'),
    @('17×38=646', '38×93=3534'),
    @('75×58=4350', '62×85=5270'),
    @('24×53=1272', '47×53=2491'),
    @('34×77=2618', '51×67=3417'),
    @('24×42=1008', '61×49=2989'),
    @('37×33=1221', '29×45=1305'),
    @('14×75=1050', '19×46=874'),
    @('50×20=1000', '59×87=5133'),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
